$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow
$win.Zoom = 147
$ws.Range("A1:M6").Select()
Write-Host "done"
